$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 123; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45203
}
